$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = -1
$ws.Range("B1").Value = -1
$ws.Range("C1").Value = -1
$ws.Range("D1").Value = -1
$ws.Range("E1").Value = 1.020373821258545
